$d = $word.ActiveDocument

# The title paragraph currently reads "7. Website". Replace the leading
# "7" with "24" (the "." Website" suffix is untouched).
$r = $d.Range(0, 1)
$r.Text = "24"

# Re-seat the document's "_GoBack" bookmark (tracking the last edit
# location) at the point right after the newly-typed "24" -- this both
# splits the run there (matching Word's own behaviour when a bookmark
# sits between two runs of identical formatting) and removes the
# bookmark from its old location near the end of the document, since
# bookmark names are unique.
$bmRange = $d.Range(2, 2)
$d.Bookmarks.Add("_GoBack", $bmRange)
